$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.477000000000004
$ws.Range("C3").Value = -10.47379999999999
$ws.Range("D3").Value = -6.905499999999992
$ws.Range("C4").Value = -12.58449999999999
$ws.Range("D9").Value = -7.259099999999997
$ws.Range("B11").Value = 5.2851
$ws.Range("B12").Value = 4.941199999999994
$ws.Range("C14").Value = -12.9412
$ws.Range("B15").Value = 4.790599999999995
$ws.Range("D15").Value = -8.385499999999995
$ws.Range("D19").Value = -7.7187
$ws.Range("D20").Value = -7.522300000000002
$ws.Range("D25").Value = -7.807000000000003
$ws.Range("C26").Value = -11.6483
$ws.Range("B27").Value = 6.378300000000004
$ws.Range("D27").Value = -9.150700000000001
$ws.Range("B28").Value = 6.348600000000002
$ws.Range("D28").Value = -8.170499999999997
$ws.Range("D30").Value = -7.216700000000007
$ws.Range("B31").Value = 4.158999999999997
$ws.Range("C31").Value = -13.7212
$ws.Range("B32").Value = 6.675399999999999
$ws.Range("D32").Value = -8.622499999999997
$ws.Range("C35").Value = -12.20790000000001
$ws.Range("B36").Value = 9.266399999999997
$ws.Range("C37").Value = -13.17319999999999
$ws.Range("B38").Value = 5.528200000000003
$ws.Range("C39").Value = -12.488
$ws.Range("C40").Value = -13.526
$ws.Range("D44").Value = -7.181200000000002
$ws.Range("C45").Value = -13.5905
$ws.Range("B46").Value = 7.159899999999999
$ws.Range("D47").Value = -7.6163
$ws.Range("C52").Value = -10.9677
$ws.Range("B54").Value = 4.678499999999999
$ws.Range("B55").Value = 5.279399999999997
$ws.Range("B56").Value = 4.837200000000003
$ws.Range("C57").Value = -14.27419999999998
$ws.Range("D58").Value = -7.932899999999997
$ws.Range("D62").Value = -8.496399999999992
$ws.Range("B67").Value = 4.894299999999992
$ws.Range("B69").Value = 5.455499999999993
$ws.Range("B72").Value = 5.603300000000002
$ws.Range("B73").Value = 8.759400000000001
$ws.Range("D77").Value = -5.840899999999998
$ws.Range("D78").Value = -7.585999999999999
$ws.Range("C81").Value = -13.61869999999999
$ws.Range("B83").Value = 5.267099999999995
$ws.Range("C83").Value = -13.576
$ws.Range("D84").Value = -8.506500000000001
$ws.Range("B86").Value = 4.830100000000003
$ws.Range("D89").Value = -6.123999999999996
$ws.Range("B91").Value = 4.923699999999998
$ws.Range("D91").Value = -6.110699999999996
$ws.Range("D92").Value = -6.051099999999998
$ws.Range("B93").Value = 6.7022
$ws.Range("D96").Value = -7.575500000000008
$ws.Range("B99").Value = 4.629299999999999
$ws.Range("C100").Value = -12.93319999999999
$ws.Range("C102").Value = -13.74159999999999
$ws.Range("D102").Value = -7.727199999999997
